$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts the existing rows 33-49
# down to 34-50 and extends the sheet dimension accordingly.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly price record.
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = 45126
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = 100112044
$ws.Cells.Item(33, 7).Value = "Perejil"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 350
$ws.Cells.Item(33, 11).Value = 800
$ws.Cells.Item(33, 12).Value = 1000
$ws.Cells.Item(33, 13).Value = 857
$ws.Cells.Item(33, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 428
$ws.Cells.Item(33, 17).Value = 2
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by
# the rest of column D.
$ws.Cells.Item(33, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
